$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $c = $ws.Range($cellRef)
    $oldStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $newValue
    $c.Style = $oldStyle
}

Set-TextValue 'D2' '51.619.45'
Set-TextValue 'E2' '  -0.71%  '
Set-TextValue 'D3' '2.795.86'
Set-TextValue 'E3' '  +0.22%  '
Set-TextValue 'E4' '  -0.02%  '
Set-TextValue 'D5' '351.52'
Set-TextValue 'E5' '  -2.20%  '
Set-TextValue 'D6' '109.07'
Set-TextValue 'E6' '  -0.44%  '
Set-TextValue 'D7' '0.553'
Set-TextValue 'E7' '  -1.01%  '
Set-TextValue 'D8' '0.999'
Set-TextValue 'E8' '  -0.02%  '
Set-TextValue 'E9' '  +5.27%  '
Set-TextValue 'D10' '39.73'
Set-TextValue 'E10' '  -1.20%  '
Set-TextValue 'E11' '  +1.08%  '
Set-TextValue 'D12' '0.0835'
Set-TextValue 'E12' '  -1.84%  '
Set-TextValue 'D13' '19.99'
Set-TextValue 'E13' '  +2.50%  '
Set-TextValue 'D14' '7.79'
Set-TextValue 'E14' '  +2.77%  '
Set-TextValue 'D15' '3.230.04'
Set-TextValue 'E15' '  -0.06%  '
Set-TextValue 'D16' '2.804.13'
Set-TextValue 'E16' '  -0.54%  '
Set-TextValue 'D17' '0.934'
Set-TextValue 'E17' '  -0.74%  '
Set-TextValue 'D18' '51.634.89'
Set-TextValue 'E18' '  -0.59%  '
Set-TextValue 'D19' '7.75'
Set-TextValue 'E19' '  +3.77%  '
Set-TextValue 'D20' '3.14'
Set-TextValue 'E20' '  +1.28%  '
Set-TextValue 'D21' '13.38'
Set-TextValue 'E21' '  +1.79%  '
Set-TextValue 'D22' '0.0₃0969'
Set-TextValue 'E22' '  -0.69%  '
Set-TextValue 'D23' '70.50'
Set-TextValue 'E23' '  +0.22%  '
Set-TextValue 'D24' '267.34'
Set-TextValue 'E24' '  -1.35%  '
Set-TextValue 'D25' '2.75'
Set-TextValue 'E25' '  -0.33%  '
Set-TextValue 'E26' '  +0.09%  '
Set-TextValue 'D27' '25.93'
Set-TextValue 'E27' '  -2.24%  '
Set-TextValue 'D28' '0.165'
Set-TextValue 'E28' '  +2.48%  '
Set-TextValue 'E29' '  -0.16%  '
Set-TextValue 'D30' '37.38'
Set-TextValue 'E30' '  +8.93%  '
Set-TextValue 'D32' '6.26'
Set-TextValue 'E32' '  +8.87%  '
Set-TextValue 'D33' '51.99'
Set-TextValue 'E33' '  +0.18%  '
Set-TextValue 'D34' '5.70'
Set-TextValue 'E34' '  +9.12%  '
Set-TextValue 'D35' '0.0445'
Set-TextValue 'E35' '  -5.85%  '
Set-TextValue 'D36' '0.0854'
Set-TextValue 'E36' '  +0.75%  '
Set-TextValue 'E37' '  -0.09%  '
Set-TextValue 'D38' '18.61'
Set-TextValue 'E38' '  -2.40%  '
Set-TextValue 'E39' '  -2.57%  '
Set-TextValue 'E40' '  -0.89%  '
Set-TextValue 'E41' '  -0.43%  '
Set-TextValue 'E42' '  -5.84%  '
Set-TextValue 'D43' '120.05'
Set-TextValue 'E43' '  +0.47%  '
Set-TextValue 'B44' 'EnergySwap'
Set-TextValue 'C44' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D44' '21.95'
Set-TextValue 'E44' '  +0.77%  '
Set-TextValue 'B45' 'WEMIXToken'
Set-TextValue 'C45' 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue 'D45' '2.19'
Set-TextValue 'E45' '  -2.66%  '
Set-TextValue 'D46' '2.135.68'
Set-TextValue 'E46' '  +2.57%  '
Set-TextValue 'D47' '3.37'
Set-TextValue 'E47' '  +3.22%  '
Set-TextValue 'E48' '  +5.87%  '
Set-TextValue 'D49' '0.226'
Set-TextValue 'E49' '  +18.18%  '
Set-TextValue 'D50' '0.909'
Set-TextValue 'E50' '  -4.22%  '
Set-TextValue 'D51' '1.36'
Set-TextValue 'E51' '  +9.37%  '
